$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "333÷2=166, 1" "637÷5=127, 2"
Replace-Text "946÷8=118, 2" "113÷6=18, 5"
Replace-Text "288÷3=96, 0" "686÷9=76, 2"
Replace-Text "205÷5=41, 0" "533÷9=59, 2"
Replace-Text "490÷9=54, 4" "628÷2=314, 0"
Replace-Text "313÷5=62, 3" "796÷3=265, 1"
Replace-Text "846÷3=282, 0" "118÷3=39, 1"
Replace-Text "127÷5=25, 2" "145÷7=20, 5"
Replace-Text "913÷5=182, 3" "492÷2=246, 0"
Replace-Text "124÷3=41, 1" "300÷2=150, 0"
Replace-Text "623÷6=103, 5" "689÷6=114, 5"
Replace-Text "571÷5=114, 1" "285÷7=40, 5"
Replace-Text "642÷6=107, 0" "586÷7=83, 5"
Replace-Text "779÷5=155, 4" "763÷2=381, 1"
Replace-Text "506÷3=168, 2" "113÷7=16, 1"
Replace-Text "647÷6=107, 5" "757÷8=94, 5"
Replace-Text "261÷3=87, 0" "631÷3=210, 1"
Replace-Text "477÷9=53, 0" "450÷6=75, 0"
Replace-Text "462÷5=92, 2" "144÷3=48, 0"
Replace-Text "464÷7=66, 2" "779÷2=389, 1"
Replace-Text "280÷4=70, 0" "201÷6=33, 3"
Replace-Text "168÷4=42, 0" "297÷4=74, 1"
Replace-Text "177÷4=44, 1" "676÷3=225, 1"
Replace-Text "868÷2=434, 0" "342÷2=171, 0"
Replace-Text "449÷3=149, 2" "732÷9=81, 3"

Write-Output "Done"
